# Update live crypto price/volume snapshot values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'71.510.63"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "'3.814.28"
$ws.Range("E3").Value = "  -0.53%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'704.51"
$ws.Range("E5").Value = "  -2.08%  "

$ws.Range("D6").Value = "'171.57"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("D7").Value = "'3.816.23"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("E11").Value = "  +5.85%  "

$ws.Range("D12").Value = "'7.44"
$ws.Range("E12").Value = "  +2.20%  "

$ws.Range("E13").Value = "  -2.14%  "

$ws.Range("D14").Value = "'36.51"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "'4.454.61"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").Value = "'3.811.07"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").Value = "'71.649.10"
$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("D18").Value = "'7.22"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").Value = "'17.53"
$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "'515.68"
$ws.Range("E21").Value = "  +4.10%  "

$ws.Range("D22").Value = "'10.49"
$ws.Range("E22").Value = "  -2.01%  "

$ws.Range("E23").Value = "  -1.22%  "

$ws.Range("D24").Value = "'84.07"
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("E25").Value = "  -3.14%  "

$ws.Range("D26").Value = "'12.59"
$ws.Range("E26").Value = "  +3.53%  "

$ws.Range("D27").Value = "'3.954.89"
$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("E28").Value = "  -2.88%  "

$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  -3.65%  "

$ws.Range("E31").Value = "  -6.95%  "

$ws.Range("D32").Value = "'7.44"
$ws.Range("E32").Value = "  -1.13%  "

$ws.Range("E33").Value = "  +0.47%  "

$ws.Range("D34").Value = "'29.24"
$ws.Range("E34").Value = "  -0.63%  "

$ws.Range("D35").Value = "'0.172"
$ws.Range("E35").Value = "  -3.39%  "

$ws.Range("D36").Value = "'9.33"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("D37").Value = "'3.770.88"
$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("D38").Value = "'0.997"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("D39").Value = "'6.54"
$ws.Range("E39").Value = "  +8.66%  "

$ws.Range("E40").Value = "  -1.53%  "

$ws.Range("D41").Value = "'2.45"
$ws.Range("E41").Value = "  +8.44%  "

$ws.Range("E42").Value = "  -1.91%  "

$ws.Range("E43").Value = "  -3.40%  "

$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "'167.01"
$ws.Range("E46").Value = "  +1.92%  "

$ws.Range("D47").Value = "'50.30"
$ws.Range("E47").Value = "  +3.12%  "

$ws.Range("E48").Value = "  -3.61%  "

$ws.Range("D49").Value = "'429.61"
$ws.Range("E49").Value = "  +2.28%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'30.40"
$ws.Range("E51").Value = "  +8.73%  "
